# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G (header "K") values are recalculated; write the new values row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 2
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 3
    13 = 0
    14 = 1
    15 = 2
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 2
    31 = 1
    32 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    40 = 0
    41 = 0
    42 = 1
    43 = 2
    44 = 0
    45 = 2
    46 = 0
    47 = 1
    48 = 1
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 1
    55 = 1
    56 = 0
    57 = 1
    58 = 0
    59 = 2
    60 = 0
    61 = 0
    62 = 1
    63 = 3
    64 = 1
    65 = 0
    66 = 2
    67 = 1
    68 = 1
    69 = 1
    71 = 0
    72 = 1
    73 = 1
    74 = 0
    75 = 0
    76 = 1
    77 = 0
    78 = 2
    79 = 0
    80 = 0
    81 = 0
    82 = 3
    83 = 1
    84 = 2
    85 = 1
    86 = 0
    87 = 0
    88 = 2
    89 = 1
    90 = 0
    91 = 1
    92 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
